# Update 2520 workload doc:
#  - Move the "9.2 preparation" row so it appears immediately before the
#    "P04 Paper 2 Article Selection" row (instead of immediately after it).
#  - Rename "P05 Paper 2 Disagreement Summary" -> "P05 Paper 2 Identify the
#    Criticism".
#  - Insert a new "10.1 preparation" row immediately before the existing
#    "10.2 preparation" row.
#
# Note: row/range objects returned by this host are positional anchors, so
# after a structural edit (Rows.Add / Row.Delete) any previously-captured
# row reference can point at a different row. To stay safe we look rows up
# by index immediately before each mutation instead of holding on to row
# objects across mutations.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Find-RowIndex($table, $prefix) {
    for ($i = 1; $i -le $table.Rows.Count; $i++) {
        $cellText = $table.Rows.Item($i).Cells.Item(1).Range.Text
        if ($cellText.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Move "9.2 preparation" to sit before "P04 Paper 2 Article Selection"
# ---------------------------------------------------------------------
$p04Index = Find-RowIndex $t "P04 Paper 2 Article Selection"

# Rows.Add(beforeRow) inserts the new row AT beforeRow's index, pushing the
# old content (and everything after it) down by one.
$newRow = $t.Rows.Add($t.Rows.Item($p04Index))
$newRow.Cells.Item(1).Range.Text = "9.2 preparation"

# The original "9.2 preparation" row (which used to sit right after P04)
# is now two rows further down: p04Index (new "9.2") -> p04Index+1 (P04) ->
# p04Index+2 (old "9.2"). Confirm before deleting so we never remove the
# wrong row if the layout ever differs from what we expect.
$oldNineTwoIndex = $p04Index + 2
$oldNineTwoText = $t.Rows.Item($oldNineTwoIndex).Cells.Item(1).Range.Text
if ($oldNineTwoText.StartsWith("9.2 preparation")) {
    $t.Rows.Item($oldNineTwoIndex).Delete()
} else {
    Write-Output ("WARNING: expected old '9.2 preparation' row at index " + $oldNineTwoIndex + " but found [" + $oldNineTwoText + "]")
}

# ---------------------------------------------------------------------
# 2) Rename P05 row text
# ---------------------------------------------------------------------
$d.Content.Find.Execute("P05 Paper 2 Disagreement Summary", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "P05 Paper 2 Identify the Criticism", 2)

# ---------------------------------------------------------------------
# 3) Insert a new "10.1 preparation" row before "10.2 preparation"
# ---------------------------------------------------------------------
$tenTwoIndex = Find-RowIndex $t "10.2 preparation"

$newTenOneRow = $t.Rows.Add($t.Rows.Item($tenTwoIndex))
$newTenOneRow.Cells.Item(1).Range.Text = "10.1 preparation"
